$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3378.8845
$ws.Range("I138").Value = 1577.0968
$ws.Range("J138").Value = 6038.6665
$ws.Range("K138").Value = 4731.2904
$ws.Range("L138").Value = 18115.9995
$ws.Range("M138").Value = 408.7096000000001
$ws.Range("N138").Value = -28395.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1490.9565
$ws.Range("I45").Value = 1020.82355
$ws.Range("K45").Value = 1020.82355
$ws.Range("M45").Value = -643.82355

$ws.Range("H61").Value = 4526.8823
$ws.Range("I61").Value = 1797.6666
$ws.Range("J61").Value = 5111.7144
$ws.Range("K61").Value = 1797.6666
$ws.Range("L61").Value = 5111.7144
$ws.Range("M61").Value = -1585.6666
$ws.Range("N61").Value = -5535.7144

$ws.Range("H74").Value = 613.44116
$ws.Range("I74").Value = 475.70834
$ws.Range("J74").Value = 944
$ws.Range("K74").Value = 475.70834
$ws.Range("L74").Value = 944
$ws.Range("M74").Value = 398.29166
$ws.Range("N74").Value = -2692

$ws.Range("H77").Value = 613.44116
$ws.Range("I77").Value = 475.70834
$ws.Range("J77").Value = 944
$ws.Range("K77").Value = 2378.5417
$ws.Range("L77").Value = 4720
$ws.Range("M77").Value = 1989.4583
$ws.Range("N77").Value = -13456

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("N89").ClearContents()

$ws.Range("H110").Value = 1497.7273
$ws.Range("J110").Value = 5087.5
$ws.Range("L110").Value = 5087.5
$ws.Range("N110").Value = -9177.5

$ws.Range("H122").Value = 3290.9048
$ws.Range("J122").Value = 5228.4287
$ws.Range("L122").Value = 15685.2861
$ws.Range("N122").Value = -20585.2861

$ws.Range("H136").Value = 4526.8823
$ws.Range("I136").Value = 1797.6666
$ws.Range("J136").Value = 5111.7144
$ws.Range("K136").Value = 5392.9998
$ws.Range("L136").Value = 15335.1432
$ws.Range("M136").Value = -2842.9998
$ws.Range("N136").Value = -20435.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1646.6842
$ws.Range("I94").Value = 1480.4667
$ws.Range("J94").Value = 2270
$ws.Range("K94").Value = 1480.4667
$ws.Range("L94").Value = 2270
$ws.Range("M94").Value = -1029.4667
$ws.Range("N94").Value = -3172

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 33000
$ws.Range("J87").Value = 33000
$ws.Range("L87").Value = 33000
$ws.Range("N87").Value = -35372

$ws.Range("H90").Value = 33000
$ws.Range("J90").Value = 33000
$ws.Range("L90").Value = 99000
$ws.Range("N90").Value = -110856

$ws.Range("H99").Value = 2094.6843
$ws.Range("I99").Value = 1450
$ws.Range("J99").Value = 2170.5293
$ws.Range("K99").Value = 1450
$ws.Range("L99").Value = 2170.5293
$ws.Range("M99").Value = 48
$ws.Range("N99").Value = -5166.5293

$ws.Range("H126").Value = 2094.6843
$ws.Range("I126").Value = 1450
$ws.Range("J126").Value = 2170.5293
$ws.Range("K126").Value = 4350
$ws.Range("L126").Value = 6511.5879
$ws.Range("M126").Value = -1880
$ws.Range("N126").Value = -11451.5879

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 5160.2188
$ws.Range("I34").Value = 184
$ws.Range("J34").Value = 5491.967
$ws.Range("K34").Value = 552
$ws.Range("L34").Value = 16475.901
$ws.Range("M34").Value = -468
$ws.Range("N34").Value = -16643.901

$ws.Range("H64").Value = 10124.667
$ws.Range("I64").Value = 1233.3334
$ws.Range("J64").Value = 13088.444
$ws.Range("K64").Value = 3700.0002
$ws.Range("L64").Value = 39265.33199999999
$ws.Range("M64").Value = -3430.0002
$ws.Range("N64").Value = -39805.33199999999

$ws.Range("H67").Value = 10124.667
$ws.Range("I67").Value = 1233.3334
$ws.Range("J67").Value = 13088.444
$ws.Range("K67").Value = 3700.0002
$ws.Range("L67").Value = 39265.33199999999
$ws.Range("M67").Value = -2764.0002
$ws.Range("N67").Value = -41137.33199999999

$ws.Range("H68").Value = 2187
$ws.Range("I68").Value = 793.3333
$ws.Range("J68").Value = 2605.1
$ws.Range("K68").Value = 2379.9999
$ws.Range("L68").Value = 7815.299999999999
$ws.Range("M68").Value = -1568.9999
$ws.Range("N68").Value = -9437.299999999999

$ws.Range("H69").Value = 11566.8
$ws.Range("I69").Value = 849.5
$ws.Range("K69").Value = 2548.5
$ws.Range("M69").Value = -1737.5

$ws.Range("H71").Value = 2187
$ws.Range("I71").Value = 793.3333
$ws.Range("J71").Value = 2605.1
$ws.Range("K71").Value = 7139.9997
$ws.Range("L71").Value = 23445.9
$ws.Range("M71").Value = -3083.9997
$ws.Range("N71").Value = -31557.9

$ws.Range("H72").Value = 11566.8
$ws.Range("I72").Value = 849.5
$ws.Range("K72").Value = 7645.5
$ws.Range("M72").Value = -3589.5

$ws.Range("H107").Value = 1019.6964
$ws.Range("I107").Value = 620.375
$ws.Range("K107").Value = 1861.125
$ws.Range("M107").Value = 58.875

$ws.Range("H132").Value = 3385.3076
$ws.Range("I132").Value = 1250.5
$ws.Range("J132").Value = 6801
$ws.Range("K132").Value = 11254.5
$ws.Range("L132").Value = 61209
$ws.Range("M132").Value = -8724.5
$ws.Range("N132").Value = -66269

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 36007
$ws.Range("J23").Value = 36007
$ws.Range("L23").Value = 36007
$ws.Range("N23").Value = -36453

$ws.Range("H70").Value = 4164.5
$ws.Range("I70").Value = 3786.182
$ws.Range("J70").Value = 5551.6665
$ws.Range("K70").Value = 3786.182
$ws.Range("L70").Value = 5551.6665
$ws.Range("M70").Value = -3516.182
$ws.Range("N70").Value = -6091.6665

$ws.Range("H73").Value = 4164.5
$ws.Range("I73").Value = 3786.182
$ws.Range("J73").Value = 5551.6665
$ws.Range("K73").Value = 3786.182
$ws.Range("L73").Value = 5551.6665
$ws.Range("M73").Value = -2850.182
$ws.Range("N73").Value = -7423.6665

$ws.Range("H74").Value = 22131
$ws.Range("J74").Value = 22131
$ws.Range("L74").Value = 22131
$ws.Range("N74").Value = -24003

$ws.Range("H77").Value = 22131
$ws.Range("J77").Value = 22131
$ws.Range("L77").Value = 66393
$ws.Range("N77").Value = -75753

$ws.Range("H102").Value = 2573.5557
$ws.Range("I102").Value = 1665.8182
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 1665.8182
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = -43.81819999999993
$ws.Range("N102").Value = -7244

$ws.Range("H113").Value = 3682.2
$ws.Range("I113").Value = 1853.6666
$ws.Range("J113").Value = 6425
$ws.Range("K113").Value = 1853.6666
$ws.Range("L113").Value = 6425
$ws.Range("M113").Value = 316.3334
$ws.Range("N113").Value = -10765

$ws.Range("H126").Value = 3469.2856
$ws.Range("I126").Value = 1833.1666
$ws.Range("J126").Value = 4123.7334
$ws.Range("K126").Value = 5499.4998
$ws.Range("L126").Value = 12371.2002
$ws.Range("M126").Value = -3029.4998
$ws.Range("N126").Value = -17311.2002

$ws.Range("H132").Value = 26318888
$ws.Range("J132").Value = 4763.353
$ws.Range("L132").Value = 14290.059
$ws.Range("N132").Value = -19350.059

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H74").Value = 29217
$ws.Range("J74").Value = 29217
$ws.Range("L74").Value = 29217
$ws.Range("N74").Value = -31213

$ws.Range("H77").Value = 29217
$ws.Range("J77").Value = 29217
$ws.Range("L77").Value = 87651
$ws.Range("N77").Value = -97635

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1875.0625
$ws.Range("I126").Value = 1285.8096
$ws.Range("K126").Value = 3857.4288
$ws.Range("M126").Value = -1387.4288
